$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"-0.2100097071880278"
$ws.Range("E2").Value = [double]"-0.1397562947531871"
$ws.Range("F2").Value = [double]"-0.1350151204931508"
$ws.Range("G2").Value = [double]"-0.2423937982455314"
$ws.Range("H2").Value = [double]"-0.1237654066428298"
$ws.Range("I2").Value = [double]"-0.2139266899315793"
$ws.Range("J2").Value = [double]"0.04865741612991594"
$ws.Range("K2").Value = [double]"0.109217188323338"
$ws.Range("L2").Value = [double]"-0.2641546473894593"
$ws.Range("M2").Value = [double]"0.01962844121283543"
$ws.Range("D3").Value = [double]"-0.01509990977053945"
$ws.Range("E3").Value = [double]"0.08162223061345381"
$ws.Range("F3").Value = [double]"0.08730196051608709"
$ws.Range("G3").Value = [double]"0.1770781441752131"
$ws.Range("H3").Value = [double]"0.1997201598806487"
$ws.Range("I3").Value = [double]"-0.0873917012978082"
$ws.Range("J3").Value = [double]"-0.06407932926335395"
$ws.Range("K3").Value = [double]"-0.04380620175755993"
$ws.Range("L3").Value = [double]"0.1332748837988889"
$ws.Range("M3").Value = [double]"0.1533272690361234"
$ws.Range("B4").Value = [double]"-0.2100097071880278"
$ws.Range("C4").Value = [double]"-0.01509990977053945"
$ws.Range("E4").Value = [double]"0.6015365415393672"
$ws.Range("F4").Value = [double]"0.622839780508571"
$ws.Range("G4").Value = [double]"0.6052692578710602"
$ws.Range("H4").Value = [double]"0.3043785453484635"
$ws.Range("I4").Value = [double]"0.9089834684659657"
$ws.Range("J4").Value = [double]"0.5539968224811245"
$ws.Range("K4").Value = [double]"0.52576023735489"
$ws.Range("L4").Value = [double]"0.5646147259402614"
$ws.Range("M4").Value = [double]"0.2943428061408312"
$ws.Range("N4").Value = [double]"-2.920454072657055e-15"
$ws.Range("O4").Value = [double]"-0.2044140082319981"
$ws.Range("R4").Value = [double]"-0.2806621929417707"
$ws.Range("S4").Value = [double]"0.563808953974267"
$ws.Range("T4").Value = [double]"0.4401659758900596"
$ws.Range("U4").Value = [double]"0.3157102365176738"
$ws.Range("B5").Value = [double]"-0.1397562947531871"
$ws.Range("C5").Value = [double]"0.08162223061345381"
$ws.Range("D5").Value = [double]"0.6015365415393672"
$ws.Range("F5").Value = [double]"0.9900306049985805"
$ws.Range("G5").Value = [double]"0.1257176556239894"
$ws.Range("H5").Value = [double]"0.8177480993505468"
$ws.Range("I5").Value = [double]"0.7676081924089456"
$ws.Range("J5").Value = [double]"0.8777559965761299"
$ws.Range("K5").Value = [double]"0.7745478205579475"
$ws.Range("L5").Value = [double]"0.2102515205998472"
$ws.Range("M5").Value = [double]"0.607083010218251"
$ws.Range("N5").Value = [double]"3.336680005718555e-15"
$ws.Range("O5").Value = [double]"-0.4942361054572194"
$ws.Range("R5").Value = [double]"-0.042330879022222"
$ws.Range("S5").Value = [double]"-0.1172357904798116"
$ws.Range("T5").Value = [double]"0.3925835194979616"
$ws.Range("U5").Value = [double]"-0.3500317328095502"
$ws.Range("B6").Value = [double]"-0.1350151204931508"
$ws.Range("C6").Value = [double]"0.08730196051608709"
$ws.Range("D6").Value = [double]"0.622839780508571"
$ws.Range("E6").Value = [double]"0.9900306049985805"
$ws.Range("G6").Value = [double]"0.1369577931509407"
$ws.Range("H6").Value = [double]"0.8322900178196406"
$ws.Range("I6").Value = [double]"0.7825289384836726"
$ws.Range("J6").Value = [double]"0.8959704693059219"
$ws.Range("K6").Value = [double]"0.8039190036363218"
$ws.Range("L6").Value = [double]"0.2290172208916044"
$ws.Range("M6").Value = [double]"0.6212141847421688"
$ws.Range("N6").Value = [double]"7.666083750325158e-16"
$ws.Range("O6").Value = [double]"-0.487397570104229"
$ws.Range("R6").Value = [double]"-0.03211179235540852"
$ws.Range("S6").Value = [double]"-0.1002609806753609"
$ws.Range("T6").Value = [double]"0.3640726428403179"
$ws.Range("U6").Value = [double]"-0.33213038900801"
$ws.Range("B7").Value = [double]"-0.2423937982455314"
$ws.Range("C7").Value = [double]"0.1770781441752131"
$ws.Range("D7").Value = [double]"0.6052692578710602"
$ws.Range("E7").Value = [double]"0.1257176556239894"
$ws.Range("F7").Value = [double]"0.1369577931509407"
$ws.Range("H7").Value = [double]"-0.0827955260901168"
$ws.Range("I7").Value = [double]"0.4895566543418153"
$ws.Range("J7").Value = [double]"0.006393420111500975"
$ws.Range("K7").Value = [double]"0.01793778858604152"
$ws.Range("L7").Value = [double]"0.9586818989706541"
$ws.Range("M7").Value = [double]"-0.1173559079185561"
$ws.Range("N7").Value = [double]"2.041337486133197e-16"
$ws.Range("O7").Value = [double]"0.09051725814192106"
$ws.Range("R7").Value = [double]"-0.0687034405025222"
$ws.Range("S7").Value = [double]"0.6193672897256314"
$ws.Range("T7").Value = [double]"0.2438475521930638"
$ws.Range("U7").Value = [double]"0.3593642545909256"
$ws.Range("B8").Value = [double]"-0.1237654066428298"
$ws.Range("C8").Value = [double]"0.1997201598806487"
$ws.Range("D8").Value = [double]"0.3043785453484635"
$ws.Range("E8").Value = [double]"0.8177480993505468"
$ws.Range("F8").Value = [double]"0.8322900178196406"
$ws.Range("G8").Value = [double]"-0.0827955260901168"
$ws.Range("I8").Value = [double]"0.5183839700041073"
$ws.Range("J8").Value = [double]"0.8143982873320751"
$ws.Range("K8").Value = [double]"0.7161035345579564"
$ws.Range("L8").Value = [double]"-0.01930103596191794"
$ws.Range("M8").Value = [double]"0.778931576145244"
$ws.Range("N8").Value = [double]"-6.631039802115322e-16"
$ws.Range("O8").Value = [double]"-0.5516707609764604"
$ws.Range("R8").Value = [double]"0.1490765808683637"
$ws.Range("S8").Value = [double]"-0.3819145055213614"
$ws.Range("T8").Value = [double]"0.1864068206580604"
$ws.Range("U8").Value = [double]"-0.5215930626895318"
$ws.Range("B9").Value = [double]"-0.2139266899315793"
$ws.Range("C9").Value = [double]"-0.0873917012978082"
$ws.Range("D9").Value = [double]"0.9089834684659657"
$ws.Range("E9").Value = [double]"0.7676081924089456"
$ws.Range("F9").Value = [double]"0.7825289384836726"
$ws.Range("G9").Value = [double]"0.4895566543418153"
$ws.Range("H9").Value = [double]"0.5183839700041073"
$ws.Range("J9").Value = [double]"0.7292238056940705"
$ws.Range("K9").Value = [double]"0.6581600365985067"
$ws.Range("L9").Value = [double]"0.5269665846687024"
$ws.Range("M9").Value = [double]"0.4125547930126223"
$ws.Range("N9").Value = [double]"7.807518436920194e-16"
$ws.Range("O9").Value = [double]"-0.3007117501436224"
$ws.Range("R9").Value = [double]"-0.2181396823853521"
$ws.Range("S9").Value = [double]"0.4226321700395524"
$ws.Range("T9").Value = [double]"0.4603379440818752"
$ws.Range("U9").Value = [double]"0.1740207184735122"
$ws.Range("B10").Value = [double]"0.04865741612991594"
$ws.Range("C10").Value = [double]"-0.06407932926335395"
$ws.Range("D10").Value = [double]"0.5539968224811245"
$ws.Range("E10").Value = [double]"0.8777559965761299"
$ws.Range("F10").Value = [double]"0.8959704693059219"
$ws.Range("G10").Value = [double]"0.006393420111500975"
$ws.Range("H10").Value = [double]"0.8143982873320751"
$ws.Range("I10").Value = [double]"0.7292238056940705"
$ws.Range("K10").Value = [double]"0.9654431488146151"
$ws.Range("L10").Value = [double]"0.07050475755270486"
$ws.Range("M10").Value = [double]"0.7697349593525429"
$ws.Range("N10").Value = [double]"9.948215812580311e-15"
$ws.Range("O10").Value = [double]"-0.5377701809645062"
$ws.Range("R10").Value = [double]"-0.2268680819084783"
$ws.Range("S10").Value = [double]"-0.1195319569650822"
$ws.Range("T10").Value = [double]"0.2836745318532706"
$ws.Range("U10").Value = [double]"-0.2919061767525117"
$ws.Range("B11").Value = [double]"0.109217188323338"
$ws.Range("C11").Value = [double]"-0.04380620175755993"
$ws.Range("D11").Value = [double]"0.52576023735489"
$ws.Range("E11").Value = [double]"0.7745478205579475"
$ws.Range("F11").Value = [double]"0.8039190036363218"
$ws.Range("G11").Value = [double]"0.01793778858604152"
$ws.Range("H11").Value = [double]"0.7161035345579564"
$ws.Range("I11").Value = [double]"0.6581600365985067"
$ws.Range("J11").Value = [double]"0.9654431488146151"
$ws.Range("L11").Value = [double]"0.0630568220275609"
$ws.Range("M11").Value = [double]"0.7819199761477985"
$ws.Range("N11").Value = [double]"2.136580471001812e-15"
$ws.Range("O11").Value = [double]"-0.5305931793137654"
$ws.Range("R11").Value = [double]"-0.2928939929855898"
$ws.Range("S11").Value = [double]"-0.09199152183265459"
$ws.Range("T11").Value = [double]"0.2060365383140566"
$ws.Range("U11").Value = [double]"-0.2405096388005949"
$ws.Range("B12").Value = [double]"-0.2641546473894593"
$ws.Range("C12").Value = [double]"0.1332748837988889"
$ws.Range("D12").Value = [double]"0.5646147259402614"
$ws.Range("E12").Value = [double]"0.2102515205998472"
$ws.Range("F12").Value = [double]"0.2290172208916044"
$ws.Range("G12").Value = [double]"0.9586818989706541"
$ws.Range("H12").Value = [double]"-0.01930103596191794"
$ws.Range("I12").Value = [double]"0.5269665846687024"
$ws.Range("J12").Value = [double]"0.07050475755270486"
$ws.Range("K12").Value = [double]"0.0630568220275609"
$ws.Range("M12").Value = [double]"-0.1462537487898222"
$ws.Range("N12").Value = [double]"1.143879312308305e-16"
$ws.Range("O12").Value = [double]"0.1085129844232405"
$ws.Range("R12").Value = [double]"-0.01903605818527438"
$ws.Range("S12").Value = [double]"0.5471289431832419"
$ws.Range("T12").Value = [double]"0.2467166122058998"
$ws.Range("U12").Value = [double]"0.2641093728914577"
$ws.Range("B13").Value = [double]"0.01962844121283543"
$ws.Range("C13").Value = [double]"0.1533272690361234"
$ws.Range("D13").Value = [double]"0.2943428061408312"
$ws.Range("E13").Value = [double]"0.607083010218251"
$ws.Range("F13").Value = [double]"0.6212141847421688"
$ws.Range("G13").Value = [double]"-0.1173559079185561"
$ws.Range("H13").Value = [double]"0.778931576145244"
$ws.Range("I13").Value = [double]"0.4125547930126223"
$ws.Range("J13").Value = [double]"0.7697349593525429"
$ws.Range("K13").Value = [double]"0.7819199761477985"
$ws.Range("L13").Value = [double]"-0.1462537487898222"
$ws.Range("N13").Value = [double]"-3.663626190122911e-15"
$ws.Range("O13").Value = [double]"-0.4190085292392897"
$ws.Range("R13").Value = [double]"-0.1226664923885096"
$ws.Range("S13").Value = [double]"-0.2273234080895403"
$ws.Range("T13").Value = [double]"0.1391403236753642"
$ws.Range("U13").Value = [double]"-0.3333457373732241"
$ws.Range("D14").Value = [double]"-2.920454072657055e-15"
$ws.Range("E14").Value = [double]"3.336680005718555e-15"
$ws.Range("F14").Value = [double]"7.666083750325158e-16"
$ws.Range("G14").Value = [double]"2.041337486133197e-16"
$ws.Range("H14").Value = [double]"-6.631039802115322e-16"
$ws.Range("I14").Value = [double]"7.807518436920194e-16"
$ws.Range("J14").Value = [double]"9.948215812580311e-15"
$ws.Range("K14").Value = [double]"2.136580471001812e-15"
$ws.Range("L14").Value = [double]"1.143879312308305e-16"
$ws.Range("M14").Value = [double]"-3.663626190122911e-15"
$ws.Range("D15").Value = [double]"-0.2044140082319981"
$ws.Range("E15").Value = [double]"-0.4942361054572194"
$ws.Range("F15").Value = [double]"-0.487397570104229"
$ws.Range("G15").Value = [double]"0.09051725814192106"
$ws.Range("H15").Value = [double]"-0.5516707609764604"
$ws.Range("I15").Value = [double]"-0.3007117501436224"
$ws.Range("J15").Value = [double]"-0.5377701809645062"
$ws.Range("K15").Value = [double]"-0.5305931793137654"
$ws.Range("L15").Value = [double]"0.1085129844232405"
$ws.Range("M15").Value = [double]"-0.4190085292392897"
$ws.Range("D18").Value = [double]"-0.2806621929417707"
$ws.Range("E18").Value = [double]"-0.042330879022222"
$ws.Range("F18").Value = [double]"-0.03211179235540852"
$ws.Range("G18").Value = [double]"-0.0687034405025222"
$ws.Range("H18").Value = [double]"0.1490765808683637"
$ws.Range("I18").Value = [double]"-0.2181396823853521"
$ws.Range("J18").Value = [double]"-0.2268680819084783"
$ws.Range("K18").Value = [double]"-0.2928939929855898"
$ws.Range("L18").Value = [double]"-0.01903605818527438"
$ws.Range("M18").Value = [double]"-0.1226664923885096"
$ws.Range("D19").Value = [double]"0.563808953974267"
$ws.Range("E19").Value = [double]"-0.1172357904798116"
$ws.Range("F19").Value = [double]"-0.1002609806753609"
$ws.Range("G19").Value = [double]"0.6193672897256314"
$ws.Range("H19").Value = [double]"-0.3819145055213614"
$ws.Range("I19").Value = [double]"0.4226321700395524"
$ws.Range("J19").Value = [double]"-0.1195319569650822"
$ws.Range("K19").Value = [double]"-0.09199152183265459"
$ws.Range("L19").Value = [double]"0.5471289431832419"
$ws.Range("M19").Value = [double]"-0.2273234080895403"
$ws.Range("D20").Value = [double]"0.4401659758900596"
$ws.Range("E20").Value = [double]"0.3925835194979616"
$ws.Range("F20").Value = [double]"0.3640726428403179"
$ws.Range("G20").Value = [double]"0.2438475521930638"
$ws.Range("H20").Value = [double]"0.1864068206580604"
$ws.Range("I20").Value = [double]"0.4603379440818752"
$ws.Range("J20").Value = [double]"0.2836745318532706"
$ws.Range("K20").Value = [double]"0.2060365383140566"
$ws.Range("L20").Value = [double]"0.2467166122058998"
$ws.Range("M20").Value = [double]"0.1391403236753642"
$ws.Range("D21").Value = [double]"0.3157102365176738"
$ws.Range("E21").Value = [double]"-0.3500317328095502"
$ws.Range("F21").Value = [double]"-0.33213038900801"
$ws.Range("G21").Value = [double]"0.3593642545909256"
$ws.Range("H21").Value = [double]"-0.5215930626895318"
$ws.Range("I21").Value = [double]"0.1740207184735122"
$ws.Range("J21").Value = [double]"-0.2919061767525117"
$ws.Range("K21").Value = [double]"-0.2405096388005949"
$ws.Range("L21").Value = [double]"0.2641093728914577"
$ws.Range("M21").Value = [double]"-0.3333457373732241"
